$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text could be misread as a number by Excel; force them to stay text
# by temporarily applying a text number format, then restoring the default "Normal" style.
$numericLooking = @(
    "D5",
    "D6",
    "D12",
    "D13",
    "D14",
    "D19",
    "D20",
    "D21",
    "D22",
    "D24",
    "D25",
    "D27",
    "D30",
    "D31",
    "D32",
    "D34",
    "D35",
    "D38",
    "D39",
    "D41",
    "D42",
    "D43",
    "D47",
    "D48",
    "D49",
    "D51"
)

foreach ($addr in $numericLooking) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply all the updated values
$updates = [ordered]@{
    "D2" = "67.523.66"
    "E2" = "  -2.22%  "
    "D3" = "2.643.18"
    "E3" = "  -3.64%  "
    "E4" = "  +0.05%  "
    "D5" = "598.28"
    "E5" = "  -0.85%  "
    "D6" = "166.84"
    "E6" = "  -1.31%  "
    "E8" = "  -0.72%  "
    "D9" = "2.642.11"
    "E9" = "  -3.66%  "
    "E10" = "  +0.67%  "
    "E11" = "  +1.30%  "
    "D12" = "0.366"
    "E12" = "  -0.44%  "
    "D13" = "5.22"
    "E13" = "  -2.44%  "
    "D14" = "28.03"
    "E14" = "  -3.06%  "
    "D15" = "3.130.99"
    "E15" = "  -3.42%  "
    "E16" = "  -3.11%  "
    "D17" = "67.373.16"
    "E17" = "  -2.45%  "
    "D18" = "2.638.18"
    "E18" = "  -3.53%  "
    "D19" = "11.82"
    "E19" = "  -0.53%  "
    "D20" = "7.88"
    "E20" = "  +1.92%  "
    "D21" = "364.82"
    "E21" = "  -2.15%  "
    "D22" = "4.41"
    "E22" = "  -3.25%  "
    "E23" = "  -3.55%  "
    "B24" = "Aptos"
    "C24" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    "D24" = "10.90"
    "E24" = "  +9.26%  "
    "B25" = "SuiNetwork"
    "C25" = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
    "D25" = "2.02"
    "E25" = "  -5.35%  "
    "E26" = "  +0.03%  "
    "D27" = "70.77"
    "E27" = "  -4.33%  "
    "D28" = "2.776.49"
    "E28" = "  -3.67%  "
    "E29" = "  -3.83%  "
    "D30" = "1.00"
    "E30" = "  +0.03%  "
    "D31" = "556.34"
    "E31" = "  -7.35%  "
    "D32" = "8.04"
    "E32" = "  -3.70%  "
    "E33" = "  -4.16%  "
    "D34" = "1.92"
    "E34" = "  -2.55%  "
    "D35" = "0.132"
    "E35" = "  -1.30%  "
    "E37" = "  -5.51%  "
    "D38" = "158.11"
    "E38" = "  -2.69%  "
    "D39" = "19.45"
    "E39" = "  -2.82%  "
    "E40" = "  -3.00%  "
    "B41" = "RenderToken"
    "C41" = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
    "D41" = "5.28"
    "E41" = "  -4.47%  "
    "B42" = "Stacks"
    "C42" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "D42" = "1.83"
    "E42" = "  -5.39%  "
    "D43" = "17.95"
    "E43" = "  -0.27%  "
    "E44" = "  -6.69%  "
    "E45" = "  +0.08%  "
    "B46" = "BabyDogeCoin"
    "C46" = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
    "D46" = "0.0₆0304"
    "E46" = "  -3.56%  "
    "B47" = "OKB"
    "C47" = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
    "D47" = "40.14"
    "E47" = "  -1.96%  "
    "D48" = "0.594"
    "E48" = "  -2.40%  "
    "D49" = "154.05"
    "E49" = "  -2.37%  "
    "E50" = "  -2.41%  "
    "D51" = "1.74"
    "E51" = "  -4.27%  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Restore default styling on the cells we force-formatted as text
foreach ($addr in $numericLooking) {
    $ws.Range($addr).Style = "Normal"
}
